# "Add files via upload" / "Definitiva dei compiti e ruoli"
# Rename sheets, finish the COMPITI (tasks) sheet and build out the new
# RUOLI (roles) sheet.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# ---------------------------------------------------------------------
# Rename the sheets
# ---------------------------------------------------------------------
$ws1.Name = "COMPITI"
$ws2.Name = "RUOLI"

# ---------------------------------------------------------------------
# COMPITI (sheet1) - finish off the task/module list
# ---------------------------------------------------------------------

# The "X" mark for the "DOC CASI D'USO" row moves from column D to column E
$ws1.Range("D2").ClearContents()
$ws1.Range("E2").Value = "X"

# "DOC DIAGR SEQ" row also gets an "X" under BUCCHERI
$ws1.Range("B4").Value = "X"

# New row: DOC SPIEGAZIONE, marked under INTINI
$ws1.Range("A5").Value = "DOC SPIEGAZIONE"
$ws1.Range("D5").Value = "X"

# Spacer row with a colored fill to separate the two tables
$ws1.Range("A8:E8").Interior.ThemeColor = 7

# "MODULI SW:" header moves from row 6 down to row 9
$ws1.Range("A6").ClearContents()
$ws1.Range("A9").Value = "MODULI SW:"

# New module rows with their owners marked
$ws1.Range("A10").Value = "DatiCondivisi"
$ws1.Range("E10").Value = "X"

$ws1.Range("A11").Value = "thVasca"
$ws1.Range("C11").Value = "X"

$ws1.Range("A12").Value = "Posizione"
$ws1.Range("D12").Value = "X"

$ws1.Range("A13").Value = "Acqua"
$ws1.Range("B13").Value = "X"

$ws1.Range("A14").Value = "Barca"
$ws1.Range("B14").Value = "X"

$ws1.Range("A15").Value = "Output"
$ws1.Range("C15").Value = "X"

$ws1.Range("A16").Value = "Input"
$ws1.Range("D16").Value = "X"

# Print setup for COMPITI
$ws1.PageSetup.PaperSize = 9
$ws1.PageSetup.Orientation = 1

# ---------------------------------------------------------------------
# RUOLI (sheet2) - brand new roles table
# ---------------------------------------------------------------------

$ws2.Columns.Item(1).ColumnWidth = 15.833333333333334
$ws2.Columns.Item(2).ColumnWidth = 10.666666666666666
$ws2.Columns.Item(3).ColumnWidth = 13.666666666666666

$ws2.Range("A1").Value = "RUOLI"
$ws2.Range("B1").Value = "BUCCHERI"
$ws2.Range("C1").Value = "GHISLANZONI"
$ws2.Range("D1").Value = "INTINI"
$ws2.Range("E1").Value = "MAFFEI"

$ws2.Range("A3").Value = "CAPO PROGETTO"
$ws2.Range("D3").Value = "X"
$ws2.Range("A3:E3").Interior.ThemeColor = 7

$ws2.Range("A5").Value = "SVILUPPATORI"
$ws2.Range("B5").Value = "X"
$ws2.Range("C5").Value = "X"
$ws2.Range("D5").Value = "X"
$ws2.Range("E5").Value = "X"
$ws2.Range("A5:E5").Interior.ThemeColor = 7

$ws2.Range("A7").Value = "TESTER"
$ws2.Range("B7").Value = "X"
$ws2.Range("C7").Value = "X"
$ws2.Range("D7").Value = "X"
$ws2.Range("E7").Value = "X"
$ws2.Range("A7:E7").Interior.ThemeColor = 7

# ---------------------------------------------------------------------
# Selections / active sheet - COMPITI ends up the visible/active tab
# ---------------------------------------------------------------------
$ws2.Activate()
$ws2.Range("D14").Select()

$ws1.Activate()
$ws1.Range("E2").Select()
